# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 2
    3  = 1
    4  = 1
    6  = 1
    7  = 0
    8  = 2
    9  = 2
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 3
    15 = 3
    16 = 0
    17 = 1
    18 = 1
    19 = 3
    20 = 1
    21 = 0
    22 = 1
    23 = 1
    24 = 1
    25 = 0
    26 = 2
    27 = 4
    28 = 0
    29 = 2
    30 = 0
    31 = 1
    32 = 0
    33 = 1
    34 = 2
    35 = 0
    36 = 0
    37 = 2
    38 = 1
    39 = 1
    40 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
